$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 2) corrections:
# B2: "unnamed: 1_level_1" -> "total"
$ws.Range("B2").Value = "total"
# F2: "unnamed: 5_level_1" -> "total"
$ws.Range("F2").Value = "total"
